# "update some camera info"
# Insert a new "Date_off" column between Date_on (E) and Tree_sp (old F, now G),
# backfill a couple of the existing Arb camera rows with their off-dates plus
# missing Camera_ID ("NA"), and append four new Ter camera-trap rows (15-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new Date_off column right after Date_on ---
$ws.Columns("F").Insert()
$ws.Range("F1").Value = "Date_off"

# --- 2. Backfill existing rows 15 & 16 (cameras that have since been taken off) ---
$ws.Range("C15").Value = "NA"
$ws.Range("D15").Value = "064"
$ws.Range("F15").Value = 44122   # 2020-10-18

$ws.Range("C16").Value = "NA"
$ws.Range("F16").Value = 44123   # 2020-10-19

# --- 3. Give the new rows' Date_on/Date_off columns the same date format as above ---
$ws.Range("E16:F16").Copy()
$ws.Range("E17:F20").PasteSpecial(-4122)

# --- 4. Append the four new terrestrial camera rows ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Ter"
$ws.Range("C17").Value = "NA"
$ws.Range("D17").Value = "053"
$ws.Range("E17").Value = 44092   # 2020-09-18
$ws.Range("G17").Value = "NA"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Ter"
$ws.Range("C18").Value = "NA"
$ws.Range("D18").Value = "025"
$ws.Range("E18").Value = 44092   # 2020-09-18
$ws.Range("F18").Value = 44258   # 2021-03-03
$ws.Range("G18").Value = "NA"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Ter"
$ws.Range("C19").Value = "NA"
$ws.Range("D19").Value = "001"
$ws.Range("E19").Value = 44092   # 2020-09-18
$ws.Range("F19").Value = 44258   # 2021-03-03
$ws.Range("G19").Value = "NA"

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Ter"
$ws.Range("C20").Value = "053"
$ws.Range("E20").Value = 44092   # 2020-09-18
$ws.Range("F20").Value = 44258   # 2021-03-03
$ws.Range("G20").Value = "NA"

# --- 5. Match the final selection left by the editor ---
$ws.Range("F16").Select() | Out-Null
